$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.931.08"
$ws.Range("E2").Value = "  -1.25%  "
$ws.Range("D3").Value = "2.619.04"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.91"
$ws.Range("E5").Value = "  +1.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.93"
$ws.Range("E6").Value = "  -0.71%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.583"
$ws.Range("E8").Value = "  -0.82%  "
$ws.Range("D9").Value = "2.619.15"
$ws.Range("E9").Value = "  +0.31%  "
$ws.Range("E10").Value = "  +1.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.60"
$ws.Range("E11").Value = "  -0.79%  "
$ws.Range("E12").Value = "  +0.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.360"
$ws.Range("E13").Value = "  +2.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.07"
$ws.Range("E14").Value = "  -0.87%  "
$ws.Range("D15").Value = "3.086.98"
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("D16").Value = "62.784.64"
$ws.Range("E16").Value = "  -1.26%  "
$ws.Range("E17").Value = "  -1.31%  "
$ws.Range("D18").Value = "2.608.36"
$ws.Range("E18").Value = "  +1.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.26"
$ws.Range("E19").Value = "  +0.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.42"
$ws.Range("E20").Value = "  +1.87%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "338.86"
$ws.Range("E21").Value = "  -0.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.79"
$ws.Range("E22").Value = "  +1.07%  "
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.49"
$ws.Range("E24").Value = "  -3.64%  "
$ws.Range("B25").Value = "SuiNetwork"
$ws.Range("C25").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.57"
$ws.Range("E25").Value = "  -0.87%  "
$ws.Range("B26").Value = "Fetch.AI"
$ws.Range("C26").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.62"
$ws.Range("E26").Value = "  -0.81%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.64"
$ws.Range("E27").Value = "  +3.86%  "
$ws.Range("E28").Value = "  -1.63%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.89"
$ws.Range("E30").Value = "  -1.09%  "
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "524.23"
$ws.Range("E31").Value = "  +6.05%  "
$ws.Range("E32").Value = "  +1.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.75"
$ws.Range("E33").Value = "  +1.26%  "
$ws.Range("D34").Value = "0.0₃0800"
$ws.Range("E34").Value = "  -1.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "170.97"
$ws.Range("E35").Value = "  -2.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.15"
$ws.Range("E36").Value = "  +13.81%  "
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.402"
$ws.Range("E38").Value = "  -0.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.98"
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.87"
$ws.Range("E40").Value = "  +8.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "170.64"
$ws.Range("E41").Value = "  +2.61%  "
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.73"
$ws.Range("E43").Value = "  -0.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.15"
$ws.Range("E44").Value = "  +2.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0567"
$ws.Range("E45").Value = "  +4.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.621"
$ws.Range("E46").Value = "  -0.65%  "
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0959"
$ws.Range("E47").Value = "  +0.25%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0239"
$ws.Range("E48").Value = "  -0.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.42"
$ws.Range("E49").Value = "  -0.18%  "
$ws.Range("E50").Value = "  +0.46%  "
$ws.Range("E51").Value = "  -1.15%  "
